$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.568.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.87%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.666.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.69%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4802"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2633"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.17%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06158"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.87%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07087"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.31%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.665.52"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.66%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.88"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.95%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.5949"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.14%  "

$ws.Range("E14").Value = "  -4.28%  "

$ws.Range("E15").Value = "  +1.88%  "

$ws.Range("E16").Value = "  +0.05%  "

$ws.Range("E17").Value = "  +0.04%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.561.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.93%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006768"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.10%  "

$ws.Range("E20").Value = "  +0.57%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.879.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.66%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.462"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.692"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.63%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.333"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "134.90"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.32%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.42%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.409"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.18%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "104.92"
$ws.Range("D28").Style = "Normal"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.690"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.24%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.961"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.14%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.676"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.22%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.07671"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04338"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.56%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9997"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.618"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.61%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6143"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.53%  "

$ws.Range("E37").Value = "  +1.22%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.601"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.89%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8614"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.22%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.000"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01510"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.75%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.876"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.66%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "97.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.87%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3775"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.89%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.710"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.52%  "

$ws.Range("E46").Value = "  +0.31%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.224"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.68%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05265"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.32%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "29.51"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.53%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.372"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.05%  "

$ws.Range("E51").Value = "  +0.19%  "
